# Updating dev shortcuts to goto github repos now
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Dev" / "Docs" shortcut labels with the GitHub repo shortcuts.
# (Alt+Enter-style manual line break within the cell text.)
$ws.Range("I2").Value = "GitHub\`nPCHenry"
$ws.Range("J2").Value = "GitHub\`nVictanya"

# The longer two-line labels need word wrap turned on so they render on
# separate lines inside the cell.
$ws.Range("I2:J2").WrapText = $true

# Update the active selection left on the sheet.
[void]$ws.Range("J3").Select()
